$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 63: 61 | 1133 | 248 ---
$ws.Range("A62").Copy()
$ws.Range("A63").PasteSpecial(-4122)
$ws.Range("A63").Value = 61
$ws.Range("B63").Value = "'1133"
$ws.Range("B63").Style = "Normal"
$ws.Range("C63").Value = 248

# --- Row 64: 62 | 1585/60 | 248 ---
$ws.Range("A62").Copy()
$ws.Range("A64").PasteSpecial(-4122)
$ws.Range("A64").Value = 62
$ws.Range("B64").Value = "1585/60"
$ws.Range("C64").Value = 248

# --- Row 65: 63 | 1064/3 | 248 ---
$ws.Range("A62").Copy()
$ws.Range("A65").PasteSpecial(-4122)
$ws.Range("A65").Value = 63
$ws.Range("B65").Value = "1064/3"
$ws.Range("C65").Value = 248
